# Refresh the "cryptos" price/volume snapshot (GitHub Actions style update).
# Most Price cells (column D) look like plain numbers to Excel, so a leading
# apostrophe is used to keep them stored as text (matching the original
# inlineStr cells) instead of letting Excel auto-convert them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.892.33"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "1.874.14"
$ws.Range("E3").Value = "  -1.07%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'0.7378"
$ws.Range("E5").Value = "  -4.82%  "
$ws.Range("D6").Value = "'242.35"
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'0.3159"
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("D9").Value = "'0.07206"
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("D10").Value = "'24.69"
$ws.Range("D11").Value = "'0.08409"
$ws.Range("E11").Value = "  -3.29%  "
$ws.Range("D12").Value = "'0.7494"
$ws.Range("E12").Value = "  -3.31%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.419"
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.854.30"
$ws.Range("E14").Value = "  -6.00%  "
$ws.Range("D15").Value = "'92.48"
$ws.Range("E15").Value = "  -2.23%  "
$ws.Range("D16").Value = "29.878.98"
$ws.Range("E16").Value = "  -1.15%  "
$ws.Range("D17").Value = "'6.093"
$ws.Range("E17").Value = "  -1.85%  "
$ws.Range("E18").Value = "  -2.73%  "
$ws.Range("D19").Value = "'243.45"
$ws.Range("E19").Value = "  -1.09%  "
$ws.Range("D20").Value = "'0.000007819"
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("D21").Value = "'0.9990"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").Value = "2.123.77"
$ws.Range("E22").Value = "  -9.80%  "
$ws.Range("D23").Value = "'7.983"
$ws.Range("E23").Value = "  -2.29%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").Value = "'0.1550"
$ws.Range("E25").Value = "  -5.14%  "
$ws.Range("D26").Value = "'9.267"
$ws.Range("E26").Value = "  -2.75%  "
$ws.Range("D27").Value = "'165.24"
$ws.Range("E27").Value = "  +1.18%  "
$ws.Range("D28").Value = "'18.59"
$ws.Range("E28").Value = "  -1.41%  "
$ws.Range("D29").Value = "'2.031"
$ws.Range("E29").Value = "  -1.08%  "
$ws.Range("D30").Value = "'1.500"
$ws.Range("E30").Value = "  +4.81%  "
$ws.Range("D31").Value = "'4.585"
$ws.Range("E31").Value = "  +1.36%  "
$ws.Range("D32").Value = "'1.531"
$ws.Range("E32").Value = "  -0.92%  "
$ws.Range("D33").Value = "'4.261"
$ws.Range("E33").Value = "  +3.02%  "
$ws.Range("D34").Value = "'0.05310"
$ws.Range("E34").Value = "  -3.27%  "
$ws.Range("D35").Value = "'1.235"
$ws.Range("E35").Value = "  -1.31%  "
$ws.Range("D36").Value = "'0.7527"
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("D37").Value = "'0.9998"
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("D38").Value = "'2.697"
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("D39").Value = "'0.01954"
$ws.Range("E39").Value = "  -0.60%  "
$ws.Range("D40").Value = "'2.754"
$ws.Range("E40").Value = "  -1.19%  "
$ws.Range("D41").Value = "'0.4515"
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("D42").Value = "1.115.76"
$ws.Range("E42").Value = "  +0.93%  "
$ws.Range("D43").Value = "'6.038"
$ws.Range("E43").Value = "  -0.69%  "
$ws.Range("D44").Value = "'72.26"
$ws.Range("E44").Value = "  -2.75%  "
$ws.Range("D45").Value = "'0.8566"
$ws.Range("E45").Value = "  +0.43%  "
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").Value = "'103.41"
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("B48").Value = "SynthetixNetwork"
$ws.Range("C48").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D48").Value = "'3.104"
$ws.Range("E48").Value = "  +3.68%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "'7.644"
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("D50").Value = "'1.836"
$ws.Range("E50").Value = "  -2.76%  "
$ws.Range("D51").Value = "2.021.88"
$ws.Range("E51").Value = "  -7.05%  "
